$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.652.91"
$ws.Range("E2").Value = "  +2.07%  "
$ws.Range("D3").Value = "3.751.23"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "601.31"
$ws.Range("E5").Value = "  +1.61%  "
$ws.Range("D6").Value = "168.62"
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("D7").Value = "3.748.18"
$ws.Range("E7").Value = "  +1.40%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.70%  "
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "6.30"
$ws.Range("E11").Value = "  +2.63%  "
$ws.Range("D12").Value = "0.462"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "38.19"
$ws.Range("E13").Value = "  +1.34%  "
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").Value = "4.377.96"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "3.762.85"
$ws.Range("E16").Value = "  +1.68%  "
$ws.Range("D17").Value = "68.715.90"
$ws.Range("E17").Value = "  +2.13%  "
$ws.Range("D18").Value = "7.27"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("E19").Value = "  +0.51%  "
$ws.Range("D20").Value = "17.09"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("E21").Value = "  +19.47%  "
$ws.Range("D22").Value = "495.19"
$ws.Range("E22").Value = "  +1.47%  "
$ws.Range("D23").Value = "0.729"
$ws.Range("E23").Value = "  +1.23%  "
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").Value = "  +6.64%  "
$ws.Range("D25").Value = "85.49"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +0.57%  "
$ws.Range("D27").Value = "12.33"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("D28").Value = "10.20"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("D30").Value = "2.52"
$ws.Range("E30").Value = "  +7.17%  "
$ws.Range("D31").Value = "2.98"
$ws.Range("E31").Value = "  +2.25%  "
$ws.Range("D32").Value = "7.93"
$ws.Range("E32").Value = "  +3.13%  "
$ws.Range("D33").Value = "31.98"
$ws.Range("E33").Value = "  +1.19%  "
$ws.Range("D34").Value = "3.898.50"
$ws.Range("E34").Value = "  +1.45%  "
$ws.Range("D35").Value = "3.687.64"
$ws.Range("E35").Value = "  +1.45%  "
$ws.Range("E36").Value = "  +1.43%  "
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("E38").Value = "  +2.42%  "
$ws.Range("D39").Value = "5.85"
$ws.Range("E39").Value = "  +2.06%  "
$ws.Range("E40").Value = "  +0.44%  "
$ws.Range("D41").Value = "0.323"
$ws.Range("E41").Value = "  +0.70%  "
$ws.Range("D42").Value = "443.36"
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("D43").Value = "48.91"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +1.07%  "
$ws.Range("D45").Value = "2.87"
$ws.Range("E45").Value = "  +3.23%  "
$ws.Range("D46").Value = "8.46"
$ws.Range("E46").Value = "  +2.35%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("D48").Value = "40.34"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("D49").Value = "2.828.50"
$ws.Range("E49").Value = "  +1.85%  "
$ws.Range("D50").Value = "141.63"
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("D51").Value = "0.0355"
$ws.Range("E51").Value = "  +2.95%  "
